# Reorders the "Recorded By" (column G) list of names in each row so that
# the first-listed recorder is moved to the end of the list (left rotation).
# Examples:
#   "a, b"    -> "b, a"
#   "a, b, c" -> "b, c, a"
#   "a"       -> "a"   (unchanged, nothing to reorder)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row   # xlUp = -4162, column 7 = G

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ",\s*"

    if ($parts.Count -gt 1) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value2 = $rotated
    }
}
